# Clear tissue_free_text (column S) values for rows 6-34 on the "Tier 1_obs"
# sheet wherever the free text simply duplicated the tissue label (blood,
# kidney, lymph node, ...). Per the commit message, tissue_free_text should
# only be populated when it differs from the label, so these redundant
# entries are removed (cell content cleared, cell itself kept).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tier 1_obs")

for ($row = 6; $row -le 34; $row++) {
    $ws.Cells.Item($row, 19).Value = ""
}
